$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Food", "fae", "2023-03-08", "234.0"),
    @("Entertainment", "tset", "2023-03-08", "344.0"),
    @("Food", "halla", "2023-03-08", "123.0"),
    @("Rent", "rtyu", "2023-03-08", "5678.0")
)

$startRow = 20
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]

    # Columns A and B are plain words - safe to assign directly as text.
    $ws.Cells.Item($r, 1).Value = $rowValues[0]
    $ws.Cells.Item($r, 2).Value = $rowValues[1]

    # Columns C and D hold values that look like a date / a number
    # ("2023-03-08", "234.0"). The source workbook stores these as plain
    # shared-string text (no numeric/date conversion, no special cell
    # style). Typing them straight into .Value would make Excel coerce
    # them into a real date serial / number. Instead, build the literal
    # text via a formula and paste back only the computed value - this
    # keeps the cell as plain text (t="s") without touching cell styles
    # (no quotePrefix / numFmt needed, unlike a leading apostrophe).
    $dateCell = $ws.Cells.Item($r, 3)
    $dateCell.Formula = "=""" + $rowValues[2] + """"
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163)

    $amountCell = $ws.Cells.Item($r, 4)
    $amountCell.Formula = "=""" + $rowValues[3] + """"
    $amountCell.Copy()
    $amountCell.PasteSpecial(-4163)
}
$excel.CutCopyMode = $false
